$wb = $excel.ActiveWorkbook

# Set the "Autenticação" value on the "INFORMAÇÕES GERAIS" sheet (B6) to "No Auth"
$wsInfo = $wb.Worksheets.Item("INFORMAÇÕES GERAIS")
$wsInfo.Range("B6").Value = "No Auth"

# Update the view on "Casos de Testes": zoom out a bit and move the selection
$wsCasos = $wb.Worksheets.Item("Casos de Testes")
$wsCasos.Activate() | Out-Null
$wsCasos.Application.ActiveWindow.Zoom = 63
$wsCasos.Range("I12").Select() | Out-Null

# Move selection on "INFORMAÇÕES GERAIS" and make it the active sheet/tab again
$wsInfo.Activate() | Out-Null
$wsInfo.Range("D9").Select() | Out-Null
